# Add a new "Longest Common Subsequence" DP problem row (row 22) to Sheet1,
# mirroring the existing rows' layout (Name | Description | Solution | Link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing row (21) onto the new row (22)
# so the new cells pick up the same cell styles (Neutral / Normal / Hyperlink)
# instead of minting new style entries.
$ws.Range("A21:D21").Copy() | Out-Null
$ws.Range("A22:D22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Fill in the new row's content, in column order, so shared-string entries
# are appended in the same order as the source edit.
$ws.Range("A22").Value = "Longest Common Subsequence"
$ws.Range("B22").Value = "Return length of longest common subsequence"
$ws.Range("C22").Value = "Use 2D DP array size of input length + 1 for base case 0. Iterate over strings comparing characters. If char is same, increase by one to the i-1,j-1 value.Else pick max from left and up value."

# Hyperlink the Link cell; Hyperlinks.Add also writes the cell's display text.
$ws.Hyperlinks.Add($ws.Range("D22"), "https://leetcode.com/problems/longest-common-subsequence/") | Out-Null

# Hyperlinks.Add resets D22's style, so re-apply the Hyperlink formatting
# captured from D21.
$ws.Range("D21").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the author's final selection position.
$ws.Range("C16").Select() | Out-Null
